# Auto-generated Excel COM-interop script to apply numeric cell updates
# per the target diff. Values correspond to recalculated market-price-derived
# profit figures for various FFXIV leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4674.7144
$ws.Range("I32").Value = 1747.8572
$ws.Range("J32").Value = 6138.143
$ws.Range("K32").Value = 1747.8572
$ws.Range("L32").Value = 6138.143
$ws.Range("M32").Value = -1421.8572
$ws.Range("N32").Value = -6790.143
$ws.Range("H40").Value = 6348.875
$ws.Range("J40").Value = 5075
$ws.Range("L40").Value = 5075
$ws.Range("N40").Value = -5425
$ws.Range("H80").Value = 4258.8
$ws.Range("I80").Value = 800
$ws.Range("K80").Value = 2400
$ws.Range("M80").Value = -1402
$ws.Range("H83").Value = 4258.8
$ws.Range("I83").Value = 800
$ws.Range("K83").Value = 7200
$ws.Range("M83").Value = -2208
$ws.Range("H88").Value = 502.6
$ws.Range("J88").Value = 488
$ws.Range("L88").Value = 488
$ws.Range("N88").Value = -1300
$ws.Range("H91").Value = 502.6
$ws.Range("J91").Value = 488
$ws.Range("L91").Value = 488
$ws.Range("N91").Value = -3296
$ws.Range("H98").Value = 3351.3333
$ws.Range("I98").Value = 2749.5
$ws.Range("K98").Value = 2749.5
$ws.Range("M98").Value = -1251.5
$ws.Range("H100").Value = 2259.6
$ws.Range("I100").Value = 1844.1111
$ws.Range("J100").Value = 5999
$ws.Range("K100").Value = 1844.1111
$ws.Range("L100").Value = 5999
$ws.Range("M100").Value = -1303.1111
$ws.Range("N100").Value = -7081
$ws.Range("H101").Value = 1805.5
$ws.Range("I101").Value = 551.6
$ws.Range("J101").Value = 3895.3333
$ws.Range("K101").Value = 1654.8
$ws.Range("L101").Value = 11685.9999
$ws.Range("M101").Value = -32.80000000000018
$ws.Range("N101").Value = -14929.9999
$ws.Range("H122").Value = 3351.3333
$ws.Range("I122").Value = 2749.5
$ws.Range("K122").Value = 8248.5
$ws.Range("M122").Value = -5798.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 2179.8
$ws.Range("H110").Value = 4867.8
$ws.Range("I110").Value = 699.3333
$ws.Range("K110").Value = 699.3333
$ws.Range("M110").Value = 1345.6667
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1323.3846
$ws.Range("I20").Value = 907
$ws.Range("J20").Value = 2260.25
$ws.Range("K20").Value = 907
$ws.Range("L20").Value = 2260.25
$ws.Range("M20").Value = -660
$ws.Range("N20").Value = -2754.25
$ws.Range("H35").Value = 39950
$ws.Range("J35").Value = 39950
$ws.Range("L35").Value = 39950
$ws.Range("N35").Value = -40570
$ws.Range("H105").Value = 4157
$ws.Range("I105").Value = 4188.4
$ws.Range("K105").Value = 4188.4
$ws.Range("M105").Value = -2441.4
$ws.Range("H134").Value = 22229830
$ws.Range("I134").Value = 8592.23
$ws.Range("K134").Value = 25776.69
$ws.Range("M134").Value = -23241.69

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2373.875
$ws.Range("I16").Value = 2299.2856
$ws.Range("K16").Value = 2299.2856
$ws.Range("M16").Value = -2012.2856
$ws.Range("H31").Value = 1787.3334
$ws.Range("J31").Value = 1866
$ws.Range("L31").Value = 1866
$ws.Range("N31").Value = -2456
$ws.Range("H34").Value = 1787.3334
$ws.Range("J34").Value = 1866
$ws.Range("L34").Value = 1866
$ws.Range("N34").Value = -2270
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H58").Value = 2147.0715
$ws.Range("I58").Value = 2040.4546
$ws.Range("K58").Value = 2040.4546
$ws.Range("M58").Value = -1837.4546
$ws.Range("H62").Value = 39996
$ws.Range("J62").Value = 39996
$ws.Range("L62").Value = 39996
$ws.Range("N62").Value = -41244
$ws.Range("H65").Value = 39996
$ws.Range("J65").Value = 39996
$ws.Range("L65").Value = 199980
$ws.Range("N65").Value = -206220
$ws.Range("H113").Value = 2373.875
$ws.Range("I113").Value = 2299.2856
$ws.Range("K113").Value = 2299.2856
$ws.Range("M113").Value = -129.2856000000002
$ws.Range("H114").Value = 44946
$ws.Range("J114").Value = 44946
$ws.Range("L114").Value = 44946
$ws.Range("N114").Value = -53624
$ws.Range("H115").Value = 34999.5
$ws.Range("J115").Value = 34999.5
$ws.Range("L115").Value = 34999.5
$ws.Range("N115").Value = -37349.5
$ws.Range("H121").Value = 24750
$ws.Range("I121").Value = 24500
$ws.Range("K121").Value = 24500
$ws.Range("M121").Value = -23190
$ws.Range("H122").Value = 1445.2222
$ws.Range("I122").Value = 1399.5
$ws.Range("J122").Value = 1450.9375
$ws.Range("K122").Value = 4198.5
$ws.Range("L122").Value = 4352.8125
$ws.Range("M122").Value = -1748.5
$ws.Range("N122").Value = -9252.8125
$ws.Range("H132").Value = 5198
$ws.Range("I132").Value = 4974.905
$ws.Range("K132").Value = 14924.715
$ws.Range("M132").Value = -12394.715
$ws.Range("H136").Value = 2147.0715
$ws.Range("I136").Value = 2040.4546
$ws.Range("K136").Value = 6121.3638
$ws.Range("M136").Value = -3571.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1250
$ws.Range("J59").Value = 1250
$ws.Range("L59").Value = 3750
$ws.Range("N59").Value = -4830
$ws.Range("H124").Value = 1944.5
$ws.Range("J124").Value = 1989
$ws.Range("L124").Value = 5967
$ws.Range("N124").Value = -15787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1336.8235
$ws.Range("I102").Value = 1230.5
$ws.Range("K102").Value = 1230.5
$ws.Range("M102").Value = 391.5
$ws.Range("H122").Value = 2695.516
$ws.Range("I122").Value = 2766.1904
$ws.Range("K122").Value = 8298.5712
$ws.Range("M122").Value = -5848.5712
$ws.Range("H126").Value = 5338.364
$ws.Range("I126").Value = 7103.6665
$ws.Range("J126").Value = 3220
$ws.Range("K126").Value = 21310.9995
$ws.Range("L126").Value = 9660
$ws.Range("M126").Value = -18840.9995
$ws.Range("N126").Value = -14600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3499.2
$ws.Range("I40").Value = 3499.2
$ws.Range("K40").Value = 3499.2
$ws.Range("M40").Value = -3363.2
$ws.Range("H43").Value = 13581.429
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 13581.429
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 13581.429
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -13967.429
$ws.Range("H55").Value = 207.82353
$ws.Range("J55").Value = 249.6
$ws.Range("L55").Value = 249.6
$ws.Range("N55").Value = -595.6
$ws.Range("H136").Value = 43482176
$ws.Range("I136").Value = 3537.4211
$ws.Range("K136").Value = 10612.2633
$ws.Range("M136").Value = -8062.263300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 12625
$ws.Range("J20").Value = 11833.333
$ws.Range("L20").Value = 11833.333
$ws.Range("N20").Value = -12313.333
$ws.Range("H34").Value = 8509.666999999999
$ws.Range("J34").Value = 8509.666999999999
$ws.Range("L34").Value = 8509.666999999999
$ws.Range("N34").Value = -8915.666999999999
$ws.Range("H107").Value = 1903.7858
$ws.Range("J107").Value = 2714.7144
$ws.Range("L107").Value = 8144.1432
$ws.Range("N107").Value = -11984.1432
$ws.Range("H132").Value = 1153.5385
$ws.Range("I132").Value = 1090.5454
$ws.Range("K132").Value = 3271.6362
$ws.Range("M132").Value = -741.6361999999999

